# Weekly driver report update for 2025-04-20
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bad Drivers section: updated Critical Minutes and Good Roaming Calculation (%)
# for the Intel(R) Wi-Fi 6E AX211 160MHz - 22.250.1.2 row, plus the Totals row.
$ws.Range("C3").Value = 169
$ws.Range("D3").Value = 98.8
$ws.Range("C4").Value = 169

# Good Drivers section: fill in the Driver Vintage date for the
# Intel(R) Wi-Fi 6E AX211 160MHz - 22.150.3.1 row, keeping it as literal
# text (like the other Driver Vintage cells) instead of letting Excel
# auto-convert it to a date serial number.
$ws.Range("E12").Formula = '="2022-08-29"'
$ws.Range("E12").Copy()
$ws.Range("E12").PasteSpecial(-4163)
$excel.CutCopyMode = 0
